$d = $word.ActiveDocument

# Fetch the full WordOpenXML package for the document so we can precisely
# restructure the TECHNICAL section and register the two new paragraph
# styles (Heading3 / Heading3Char) in styles.xml in a single pass.
$full = $d.Content.WordOpenXML

# WordOpenXML always stamps every paragraph with w14:paraId/w14:textId
# attributes that are not present in the original file; strip them back out
# before we feed the XML back in so the round-trip stays clean.
$full = $full -replace ' w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+"', ''

# 1) XXX -> AAA
$full = $full.Replace("<w:t>XXX</w:t>", "<w:t>AAA</w:t>")

# 2) Replace the empty bookmark paragraph + TECHNICAL heading + "Sap developer..."
#    paragraph with the new TECHNICAL/DevOps block of paragraphs.
$oldBlock = @'
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00DF5C9E" w:rsidRDefault="00DF5C9E" w:rsidP="001A6232"><w:pPr><w:ind w:left="2160" w:firstLine="720"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>TECHNICAL</w:t></w:r></w:p><w:p w:rsidR="00C43764" w:rsidRDefault="00E16D0A" w:rsidP="00DF5C9E"><w:pPr><w:ind w:left="1440" w:firstLine="720"/></w:pPr><w:r><w:t>Sap developer,abap,pipo and ui5</w:t></w:r></w:p>
'@

$newBlock = @'
</w:p><w:p><w:pPr><w:ind w:left="2160" w:firstLine="720"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>TECHNICAL</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t>Cloud Operations Architect (DevOps) - DevOps</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Skill Details </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Cloud Computing- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - 48 months</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Shell Scripting- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - 96 months</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Python- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - 6 months</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Automation- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - 72 months</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Solution Architect- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Less than 1 year months</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Azure- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Less than 1 year months</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">AWS- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Exprience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Less than 1 year </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>monthsCompany</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Details </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>company</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> - DevOps</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>description</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> - Type: DevOps Engineer.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t>Platform: AWS Cloud, Azure Cloud.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Services: AWS EC2, RDS, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CloudFormation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Template, Lambda, Dynamo DB,</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Cloud Watch, Auto-scaling, Elastic Bean stalk, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Appdynamics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Here I manage </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tibco</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Spotfire</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> enterprise &amp; Cloud Product support. Being the only Ops member in India I got a chance to recruit &amp; build entire team of 15 members. I also worked on 4 different Projects / products simultaneously and added the hired members into these products.</w:t></w:r></w:p>
'@

if ($full.IndexOf($oldBlock) -lt 0) {
    throw "anchor block for TECHNICAL section not found"
}
$full = $full.Replace($oldBlock, $newBlock)

# 3) Register the Heading3 paragraph style right after the Normal style.
$heading3Style = @'
<w:style w:type="paragraph" w:styleId="Heading3"><w:name w:val="heading 3"/><w:basedOn w:val="Normal"/><w:link w:val="Heading3Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00285586"/><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="2"/></w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="en-IN" w:bidi="ta-IN"/></w:rPr></w:style>
'@

$normalAnchor = '<w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/><w:qFormat/></w:style>'
if ($full.IndexOf($normalAnchor) -lt 0) {
    throw "Normal style anchor not found"
}
$full = $full.Replace($normalAnchor, $normalAnchor + $heading3Style)

# 4) Register the Heading3Char character style right after the Hyperlink style.
$heading3CharStyle = @'
<w:style w:type="character" w:customStyle="1" w:styleId="Heading3Char"><w:name w:val="Heading 3 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading3"/><w:uiPriority w:val="9"/><w:rsid w:val="00285586"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="en-IN" w:bidi="ta-IN"/></w:rPr></w:style>
'@

$hyperlinkAnchor = '<w:style w:type="character" w:styleId="Hyperlink"><w:name w:val="Hyperlink"/><w:basedOn w:val="DefaultParagraphFont"/><w:uiPriority w:val="99"/><w:unhideWhenUsed/><w:rsid w:val="00DF5C9E"/><w:rPr><w:color w:val="0563C1" w:themeColor="hyperlink"/><w:u w:val="single"/></w:rPr></w:style>'
if ($full.IndexOf($hyperlinkAnchor) -lt 0) {
    throw "Hyperlink style anchor not found"
}
$full = $full.Replace($hyperlinkAnchor, $hyperlinkAnchor + $heading3CharStyle)

# Push the rebuilt package back into the document.
$d.Content.InsertXML($full)

Write-Host "Paragraph count:" $d.Paragraphs.Count
Write-Host "Styles count:" $d.Styles.Count
